# Translate the Chinese BOM annotation text to English (per commit:
# "Translated to English") and tidy up the related formatting, matching
# the author's re-save of the workbook.
#
# NOTE: the order of the .Value assignments below matters -- it controls
# the order in which new shared-string table entries are created, which
# in turn has to line up with the target workbook's string order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -- Row 17: header row for the "other tested brands" mini-table -----------
$ws.Range("A17").Value = "No.10 Other tested brand models:"
$ws.Range("B17").Value = "TI( Texas Instrument)"
$ws.Range("C17").Value = "Nexperia(Anse )"

# -- Row 18: "Commodity model" label (now wraps onto two lines) ------------
$ws.Range("B18").Value = "`nCommodity model"
$ws.Range("C18").Value = "`nCommodity model"
$ws.Rows.Item(18).RowHeight = 20.25

# -- Row 19: brand model values (unchanged text, kept for completeness) ----
$ws.Range("B19").Value = "SN74HC245PWR"
$ws.Range("C19").Value = "74HC245PW,118"

# -- Row 14: PCB thickness note ---------------------------------------------
$ws.Range("A14").Value = "PCB thickness = 1.2mm"

# -- Row 16: "No.3" note -----------------------------------------------------
$ws.Range("A16").Value = "No.3 can be removed"

# -- Selection, matching the saved cursor position in the re-uploaded file -
$ws.Range("H2").Select()
